# Apply the "new requirement separation" data update to the DATA worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers ---
# Insert three new header columns (Corequisites, Concurrent, Recommended)
# before the existing "Terms Typically Offered" column, which now moves to G1.
$ws.Range("G1").Value = $ws.Range("D1").Text
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# --- Row 2: DATA 301 ---
$ws.Range("G2").Value = $ws.Range("D2").Text
$ws.Range("C2").Value = "CPE/CSC 202; and one of the STAT 302, STAT 312, or STAT 313."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"

# --- Row 3: DATA 401 ---
$ws.Range("G3").Value = $ws.Range("D3").Text
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"

# --- Row 4: DATA 451 ---
$ws.Range("G4").Value = $ws.Range("D4").Text
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"

# --- Row 5: DATA 452 ---
$ws.Range("G5").Value = $ws.Range("D5").Text
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
